$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 13566.667
$ws.Range("L13").Value = 18900
$ws.Range("N13").Value = -19238
$ws.Range("J13").Value = 18900
$ws.Range("K34").Value = 1761.2727
$ws.Range("M34").Value = -1558.2727
$ws.Range("I34").Value = 1761.2727
$ws.Range("H34").Value = 3659.5386
$ws.Range("H36").Value = 3659.5386
$ws.Range("M36").Value = -1046.2727
$ws.Range("K36").Value = 1761.2727
$ws.Range("I36").Value = 1761.2727
$ws.Range("N129").Value = -13315.1155
$ws.Range("J129").Value = 1105.0385
$ws.Range("L129").Value = 3315.1155
$ws.Range("H129").Value = 1023.55
$ws.Range("K137").Value = 0
$ws.Range("H137").Value = 2550
$ws.Range("I137").Value = 0
$ws.Range("N137").Value = -12750
$ws.Range("M137").ClearContents()
$ws.Range("L137").Value = 7650
$ws.Range("J137").Value = 2550
$ws.Range("H138").Value = 1737.6478
$ws.Range("M138").Value = 1871.2186
$ws.Range("I138").Value = 1089.5938
$ws.Range("K138").Value = 3268.7814
$ws.Range("J138").Value = 2269.3845
$ws.Range("L138").Value = 6808.1535
$ws.Range("N138").Value = -17088.1535
$ws.Range("N141").Value = -21660.0001
$ws.Range("I141").Value = 2584.5
$ws.Range("L141").Value = 11300.0001
$ws.Range("M141").Value = -2573.5
$ws.Range("H141").Value = 3047.087
$ws.Range("J141").Value = 3766.6667
$ws.Range("K141").Value = 7753.5

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M2").Value = -449.2222
$ws.Range("L2").Value = 2406.5
$ws.Range("H2").Value = 897.5454999999999
$ws.Range("I2").Value = 562.2222
$ws.Range("K2").Value = 562.2222
$ws.Range("J2").Value = 2406.5
$ws.Range("N2").Value = -2632.5
$ws.Range("I61").Value = 2305.3
$ws.Range("L61").Value = 4559.6924
$ws.Range("K61").Value = 2305.3
$ws.Range("J61").Value = 4559.6924
$ws.Range("M61").Value = -2093.3
$ws.Range("N61").Value = -4983.6924
$ws.Range("H61").Value = 3193.394
$ws.Range("N92").Value = -37129.5
$ws.Range("L92").Value = 32137.5
$ws.Range("J92").Value = 32137.5
$ws.Range("H92").Value = 32137.5
$ws.Range("I116").Value = 562.2222
$ws.Range("M116").Value = 1731.7778
$ws.Range("N116").Value = -6994.5
$ws.Range("L116").Value = 2406.5
$ws.Range("K116").Value = 562.2222
$ws.Range("H116").Value = 897.5454999999999
$ws.Range("J116").Value = 2406.5
$ws.Range("K122").Value = 3857149.2
$ws.Range("I122").Value = 1285716.4
$ws.Range("N122").Value = -18110.5
$ws.Range("L122").Value = 13210.5
$ws.Range("H122").Value = 1072164.2
$ws.Range("M122").Value = -3854699.2
$ws.Range("J122").Value = 4403.5
$ws.Range("K132").Value = 6633
$ws.Range("H132").Value = 4351783
$ws.Range("I132").Value = 2211
$ws.Range("M132").Value = -4103
$ws.Range("J136").Value = 4559.6924
$ws.Range("N136").Value = -18779.0772
$ws.Range("L136").Value = 13679.0772
$ws.Range("H136").Value = 3193.394
$ws.Range("I136").Value = 2305.3
$ws.Range("K136").Value = 6915.900000000001
$ws.Range("M136").Value = -4365.900000000001
$ws.Range("J139").Value = 62357.5
$ws.Range("L139").Value = 62357.5
$ws.Range("H139").Value = 62357.5
$ws.Range("N139").Value = -72637.5

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("J3").Value = 2406.5
$ws.Range("K3").Value = 562.2222
$ws.Range("I3").Value = 562.2222
$ws.Range("N3").Value = -2634.5
$ws.Range("L3").Value = 2406.5
$ws.Range("M3").Value = -448.2222
$ws.Range("H3").Value = 897.5454999999999

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M16").Value = -3847534
$ws.Range("H16").Value = 2332925.5
$ws.Range("K16").Value = 3847821
$ws.Range("I16").Value = 3847821
$ws.Range("L31").Value = 10206.637
$ws.Range("I31").Value = 1594.4
$ws.Range("H31").Value = 4650.355
$ws.Range("J31").Value = 10206.637
$ws.Range("M31").Value = -1299.4
$ws.Range("N31").Value = -10796.637
$ws.Range("K31").Value = 1594.4
$ws.Range("J34").Value = 10206.637
$ws.Range("N34").Value = -10610.637
$ws.Range("L34").Value = 10206.637
$ws.Range("K34").Value = 1594.4
$ws.Range("I34").Value = 1594.4
$ws.Range("M34").Value = -1392.4
$ws.Range("H34").Value = 4650.355
$ws.Range("K58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("N58").Value = -2640.5
$ws.Range("J58").Value = 2234.5
$ws.Range("M58").ClearContents()
$ws.Range("L58").Value = 2234.5
$ws.Range("H58").Value = 2234.5
$ws.Range("M94").Value = -4343.727
$ws.Range("I94").Value = 4794.727
$ws.Range("K94").Value = 4794.727
$ws.Range("H94").Value = 4688.5356
$ws.Range("K113").Value = 3847821
$ws.Range("M113").Value = -3845651
$ws.Range("I113").Value = 3847821
$ws.Range("H113").Value = 2332925.5
$ws.Range("J136").Value = 2234.5
$ws.Range("N136").Value = -11803.5
$ws.Range("L136").Value = 6703.5
$ws.Range("H136").Value = 2234.5
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 6742.6
$ws.Range("K87").Value = 20227.8
$ws.Range("I87").Value = 6742.6
$ws.Range("M87").Value = -18979.8
$ws.Range("I90").Value = 6742.6
$ws.Range("K90").Value = 60683.4
$ws.Range("M90").Value = -54443.4
$ws.Range("H90").Value = 6742.6
$ws.Range("L114").Value = 20713.8339
$ws.Range("J114").Value = 6904.6113
$ws.Range("M114").Value = 2867.375
$ws.Range("H114").Value = 4819.769
$ws.Range("I114").Value = 128.875
$ws.Range("N114").Value = -27221.8339
$ws.Range("K114").Value = 386.625
$ws.Range("L131").Value = 3444.1464
$ws.Range("M131").Value = -19997052
$ws.Range("H131").Value = 1786741.6
$ws.Range("N131").Value = -13524.1464
$ws.Range("K131").Value = 20002092
$ws.Range("I131").Value = 6667364
$ws.Range("J131").Value = 1148.0488

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("K132").Value = 9400.049999999999
$ws.Range("N132").Value = -13034
$ws.Range("L132").Value = 7974
$ws.Range("H132").Value = 2889.878
$ws.Range("I132").Value = 3133.35
$ws.Range("J132").Value = 2658
$ws.Range("M132").Value = -6870.049999999999

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("J22").Value = 2527.2222
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("H22").Value = 2527.2222
$ws.Range("N22").Value = -3117.2222
$ws.Range("L22").Value = 2527.2222
$ws.Range("M22").ClearContents()
$ws.Range("K27").Value = 0
$ws.Range("N27").Value = -2741.2222
$ws.Range("L27").Value = 2527.2222
$ws.Range("J27").Value = 2527.2222
$ws.Range("H27").Value = 2527.2222
$ws.Range("I27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("K132").Value = 76932780
$ws.Range("N132").Value = -14284.25
$ws.Range("L132").Value = 9224.25
$ws.Range("H132").Value = 15876190
$ws.Range("I132").Value = 25644260
$ws.Range("J132").Value = 3074.75
$ws.Range("M132").Value = -76930250
$ws.Range("J136").Value = 14399.889
$ws.Range("N136").Value = -48299.667
$ws.Range("L136").Value = 43199.667
$ws.Range("H136").Value = 6681.6895
$ws.Range("I136").Value = 3208.5
$ws.Range("K136").Value = 9625.5
$ws.Range("M136").Value = -7075.5

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("K132").Value = 3115.6155
$ws.Range("N132").Value = -17510.75
$ws.Range("L132").Value = 12450.75
$ws.Range("H132").Value = 1770.7059
$ws.Range("I132").Value = 1038.5385
$ws.Range("J132").Value = 4150.25
$ws.Range("M132").Value = -585.6155000000003
$ws.Range("J136").Value = 2182.75
$ws.Range("N136").Value = -11648.25
$ws.Range("L136").Value = 6548.25
$ws.Range("H136").Value = 3425.3872
$ws.Range("I136").Value = 5684.727
$ws.Range("K136").Value = 17054.181
$ws.Range("M136").Value = -14504.181
